$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - fill previously empty cells
$ws.Range("B10").Value = "SEMP - 1"
$ws.Range("C10").Value = "SEMP TUT - 1"
$ws.Range("D10").Value = "CN"
$ws.Range("E10").Value = "AI"
$ws.Range("F10").Value = "AI"

# Row 13 - fill previously empty cells, change F13
$ws.Range("B13").Value = "SEMP - 1"
$ws.Range("C13").Value = "SEMP TUT - 1"
$ws.Range("D13").Value = "CN"
$ws.Range("E13").Value = "AI"
$ws.Range("F13").Value = "AI"

# Row 19 - fill previously empty cells, change D19/E19/F19
$ws.Range("B19").Value = "AI LAB"
$ws.Range("C19").Value = "DBMS LAB"
$ws.Range("D19").Value = "DBMS"
$ws.Range("E19").Value = "DBMS"
$ws.Range("F19").Value = "DBMS"

# Row 22 - fill previously empty cells, change D22/E22/F22
$ws.Range("B22").Value = "CN LAB"
$ws.Range("C22").Value = "CN"
$ws.Range("D22").Value = "CO"
$ws.Range("E22").Value = "CO"
$ws.Range("F22").Value = "CO"

# Row 25 - fill previously empty cells, change D25/E25/F25
$ws.Range("B25").Value = "CN"
$ws.Range("C25").Value = "AI"
$ws.Range("D25").Value = "PSE"
$ws.Range("E25").Value = "PSE"
$ws.Range("F25").Value = "PSE"
